$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text so values like "1.001" or "0.9997"
# are not auto-converted into numbers by Excel, matching the original inlineStr cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.060.56"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.830.11"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "244.55"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "0.6328"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.07535"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "0.2944"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "23.11"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").Value = "0.07704"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "1.829.72"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "5.001"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "0.6695"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "83.23"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "0.000009617"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "6.074"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "29.070.09"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "12.60"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").Value = "226.48"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "7.150"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "160.65"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +4.49%  "
$ws.Range("D26").Value = "8.520"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "17.94"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "4.149"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").Value = "4.067"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "0.05480"
$ws.Range("E31").Value = "  +5.41%  "
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "1.858"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").Value = "0.7458"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").Value = "1.138"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "2.660"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("D37").Value = "1.244.57"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").Value = "2.755"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "0.01784"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "6.628"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").Value = "0.9038"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "101.40"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "1.975.47"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("D46").Value = "65.07"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").Value = "0.5101"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").Value = "8.965"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "1.660"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "0.05787"
$ws.Range("E51").Value = "  +0.86%  "

# Restore default (unstyled) formatting now that the text values are locked in,
# so the cells have no explicit style index, matching the original workbook.
$priceRange.Style = "Normal"

